$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2
$ws.Range("B12").Value = 63
$ws.Range("C12").Value = -22
$ws.Range("E12").Value = "41/252"
